$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1953125
    "C2" = 0.546875
    "J2" = 0.0234375
    "P2" = 0.13671875
    "S2" = 0.09765625
    "B3" = 0.01360544217687075
    "C3" = 0.006802721088435374
    "J3" = 0.06122448979591837
    "P3" = 0.7210884353741497
    "S3" = 0.1972789115646258
    "J4" = 0.02702702702702703
    "P4" = 0.7297297297297297
    "S4" = 0.2432432432432433
    "J5" = 0.5
    "P5" = 0.5
    "B6" = 0.0546218487394958
    "D6" = 0.004201680672268907
    "F6" = 0.05042016806722689
    "J6" = 0.2478991596638656
    "O6" = 0.01260504201680672
    "Q6" = 0.1470588235294118
    "R6" = 0.1008403361344538
    "S6" = 0.3823529411764706
    "B7" = 0.08333333333333333
    "D7" = 0.0119047619047619
    "E7" = 0.005952380952380952
    "F7" = 0.04166666666666666
    "J7" = 0.1845238095238095
    "O7" = 0.02380952380952381
    "Q7" = 0.1011904761904762
    "R7" = 0.09523809523809523
    "S7" = 0.4523809523809524
    "B8" = 0.06511627906976744
    "D8" = 0.004651162790697674
    "F8" = 0.04883720930232558
    "J8" = 0.09767441860465116
    "O8" = 0.006976744186046512
    "Q8" = 0.2
    "R8" = 0.1116279069767442
    "S8" = 0.4651162790697674
    "B9" = 0.0896551724137931
    "D9" = 0.0103448275862069
    "F9" = 0.07241379310344828
    "J9" = 0.1310344827586207
    "O9" = 0.0103448275862069
    "Q9" = 0.1344827586206896
    "R9" = 0.1068965517241379
    "S9" = 0.4448275862068966
    "B10" = 0.08230452674897119
    "D10" = 0.02126200274348422
    "E10" = 0.001371742112482853
    "F10" = 0.08161865569272976
    "J10" = 0.1255144032921811
    "O10" = 0.0102880658436214
    "Q10" = 0.1934156378600823
    "R10" = 0.09465020576131687
    "S10" = 0.3895747599451303
    "G11" = 0.1120331950207469
    "J11" = 0.08713692946058091
    "K11" = 0.1618257261410788
    "L11" = 0.6182572614107884
    "S11" = 0.02074688796680498
    "G12" = 0.7243589743589743
    "J12" = 0.1858974358974359
    "K12" = 0.02564102564102564
    "L12" = 0.05128205128205128
    "S12" = 0.01282051282051282
    "G13" = 0.8378378378378378
    "J13" = 0.1351351351351351
    "S13" = 0.02702702702702703
    "F15" = 0.01932367149758454
    "H15" = 0.1594202898550725
    "I15" = 0.1111111111111111
    "J15" = 0.4251207729468599
    "K15" = 0.04347826086956522
    "M15" = 0.02898550724637681
    "O15" = 0.05314009661835749
    "S15" = 0.1594202898550725
    "F16" = 0.006024096385542169
    "H16" = 0.1506024096385542
    "I16" = 0.0783132530120482
    "J16" = 0.4578313253012048
    "K16" = 0.09036144578313253
    "M16" = 0.01204819277108434
    "N16" = 0.006024096385542169
    "O16" = 0.07228915662650602
    "S16" = 0.1265060240963855
    "F17" = 0.01538461538461539
    "H17" = 0.1362637362637363
    "I17" = 0.1164835164835165
    "J17" = 0.4725274725274725
    "K17" = 0.0945054945054945
    "M17" = 0.01098901098901099
    "O17" = 0.04835164835164835
    "S17" = 0.1054945054945055
    "F18" = 0.01556420233463035
    "H18" = 0.1361867704280156
    "I18" = 0.1361867704280156
    "J18" = 0.4474708171206226
    "K18" = 0.07392996108949416
    "M18" = 0.01945525291828794
    "O18" = 0.06614785992217899
    "S18" = 0.1050583657587549
    "F19" = 0.01298701298701299
    "H19" = 0.2005772005772006
    "I19" = 0.119047619047619
    "J19" = 0.4054834054834055
    "K19" = 0.08152958152958154
    "M19" = 0.0165945165945166
    "N19" = 0.0007215007215007215
    "O19" = 0.06349206349206349
    "S19" = 0.09956709956709957
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
